$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Phone-number ("Nomor Telepon", column C) values must stay TEXT (leading
# zeros). Plain `.Value = "0812..."` gets auto-coerced to a number by Excel,
# so for any phone number that already exists elsewhere in the sheet we
# relocate it with Copy/PasteSpecial (values), which preserves the original
# text cell type without touching styles. Do this BEFORE overwriting the
# other columns, and in an order where every source cell is read before it
# is itself overwritten.
# ---------------------------------------------------------------------------

# C2 <- old C6 ("081233072122")
$ws.Range("C6").Copy()
$ws.Range("C2").PasteSpecial(-4163)

# C6 <- old C10 ("085745966707")
$ws.Range("C10").Copy()
$ws.Range("C6").PasteSpecial(-4163)

# C10 <- old C12 ("087777284179")
$ws.Range("C12").Copy()
$ws.Range("C10").PasteSpecial(-4163)

# C11 <- old C13 ("083834657395")
$ws.Range("C13").Copy()
$ws.Range("C11").PasteSpecial(-4163)

$excel.CutCopyMode = 0

# C12 needs a brand-new phone number ("087654321234") that doesn't exist
# anywhere else in the sheet to copy from. Build it as a text formula in a
# scratch cell outside the used range, then paste-special just the value
# (still text) into place, then clean up the scratch cell.
$ws.Range("Z1").Formula = "=""087654321234"""
$ws.Range("Z1").Copy()
$ws.Range("C12").PasteSpecial(-4163)
$excel.CutCopyMode = 0
$ws.Range("Z1").ClearContents()

# ---------------------------------------------------------------------------
# Names (column B) and the rest of the fields are plain text/numbers, so a
# direct value assignment is fine.
# ---------------------------------------------------------------------------

# --- Row 2: Samsul Huda -> Hariyanto, ID Kota 2 -> 21 ---
$ws.Range("B2").Value = "Hariyanto"
$ws.Range("D2").Value = 21

# --- Row 6: Tolkha Hasan -> Reanaldo Revanzah Putra, ID Kota 1 -> 21 ---
$ws.Range("B6").Value = "Reanaldo Revanzah Putra"
$ws.Range("D6").Value = 21

# --- Row 10 (was Wawan Dwi Prasetyo) becomes Rizal Ferdian: ID 6 -> 8, D 4 -> 1, F 4 -> 3 ---
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Rizal Ferdian"
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 3
$ws.Range("G10").Value = 1

# --- Row 11 (was Achmad Chadil Auwfar) becomes Redika Angga Pratama: ID 7 -> 9, F 1 -> 3 ---
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Redika Angga Pratama"
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 3
$ws.Range("G11").Value = 1

# --- Row 12: newly imported employee "Amirah Rahmani" with a non-numeric (hash) ID ---
$ws.Range("A12").Value = "966a502429edc5a4e9222942eec72c57"
$ws.Range("B12").Value = "Amirah Rahmani"
$ws.Range("D12").Value = 17
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 3
$ws.Range("G12").Value = 1

# --- Row 13 no longer exists -- clear it so the used range shrinks back to A1:G12 ---
$ws.Range("A13:G13").ClearContents()
